# Weekly price update: insert a new record as row 81, pushing the existing
# rows 81-249 down to 82-250 (last row becomes 250).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(81).Insert()

$ws.Cells.Item(81, 1).Value2  = 3
$ws.Cells.Item(81, 2).Value2  = "Femacal de La Calera"
$ws.Cells.Item(81, 3).Value2  = "Coquimbo"
$ws.Cells.Item(81, 4).Value2  = 44581
$ws.Cells.Item(81, 5).Value2  = 5
$ws.Cells.Item(81, 6).Value2  = 100112039
$ws.Cells.Item(81, 7).Value2  = "Ciboulette"
$ws.Cells.Item(81, 8).Value2  = "Sin especificar"
$ws.Cells.Item(81, 9).Value2  = "Primera"
$ws.Cells.Item(81, 10).Value2 = 130
$ws.Cells.Item(81, 11).Value2 = 1500
$ws.Cells.Item(81, 12).Value2 = 1500
$ws.Cells.Item(81, 13).Value2 = 1500
$ws.Cells.Item(81, 14).Value2 = "$/docena de atados"
$ws.Cells.Item(81, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(81, 16).Value2 = 500
$ws.Cells.Item(81, 17).Value2 = 3
$ws.Cells.Item(81, 18).Value2 = "Hortaliza"
